$d = $word.ActiveDocument

# 1. Add a _GoBack bookmark at the very start of the title paragraph ("Planification du budget")
$titlePara = $d.Paragraphs(1)
$startRange = $d.Range($titlePara.Range.Start, $titlePara.Range.Start)
$d.Bookmarks.Add("_GoBack", $startRange) | Out-Null

# 2. Merge "43" + " " + "018" + " = " into a single run's text via Find/Replace
$d.Content.Find.Execute("43 018 = ", $true, $false, $false, $false, $false, $true, 1, $false, "43 018 = ", 2) | Out-Null

# 3. Remove bookmark that sits between "951," and "24", merge into "951,24"
$d.Bookmarks("_GoBack").Delete()
$d.Content.Find.Execute("951,24", $true, $false, $false, $false, $false, $true, 1, $false, "951,24", 2) | Out-Null

# Re-add _GoBack bookmark at the end (Word typically places _GoBack at last edit position)
